$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantity (F column) tweaks on existing rows ---
$ws.Range("F6").Value = 27
$ws.Range("F7").Value = 45
$ws.Range("F10").Value = 44
$ws.Range("F55").Value = 11
$ws.Range("F100").Value = 13
$ws.Range("F106").Value = 168
$ws.Range("F109").Value = 101
$ws.Range("F122").Value = 100
$ws.Range("F124").Value = 3
$ws.Range("F131").Value = 50
$ws.Range("F141").Value = 20
$ws.Range("F142").Value = 15
$ws.Range("F144").Value = 6
$ws.Range("F147").Value = 29
$ws.Range("F148").Value = 8

# --- Row 25: Buste Plastica forate f.to 22x30 - new delivery 5/6/2018, order 1100, qty 50 ---
$ws.Range("B25").Value = "5/6/2018"
$ws.Range("C25").Value = "1100"
$ws.Range("F25").Value = 50

# --- Row 32: Carta fotocopie fg A4 80 GR - new delivery date/order ---
$ws.Range("B32").Value = "5/6/2018"
$ws.Range("C32").Value = "165"

# --- Row 33: Carta Plotter - fill in new arrival date + order count ---
$ws.Range("B33").Value = 43256
$ws.Range("C33").Value = 4

# --- Row 45: Colla in Stick - new arrival + bigger order ---
$ws.Range("B45").Value = "5/6/2018"
$ws.Range("C45").Value = "7"
$ws.Range("F45").Value = 7

# --- Row 85: Mine 0,5 - new delivery ---
$ws.Range("B85").Value = "5/6/2018"
$ws.Range("C85").Value = "25"
$ws.Range("F85").Value = 15

# --- Row 103: Post-it grandi 76x76 pacchi da 12 - new delivery ---
$ws.Range("B103").Value = "5/6/2018"
$ws.Range("C103").Value = "24"
$ws.Range("F103").Value = 12

# --- Row 118: Scotch magic - new delivery ---
$ws.Range("B118").Value = "5/6/2018"
$ws.Range("C118").Value = "35"
$ws.Range("F118").Value = 6

# --- Row 119: Scotch rotolo grande - new delivery ---
$ws.Range("B119").Value = "5/6/2018"
$ws.Range("C119").Value = "14"
$ws.Range("F119").Value = 9

# --- Update selection cursor to match the author's last-saved position ---
$ws.Range("D123").Select()
